# Update countries & provincias Spain
# Refreshes the COVID case table ("Pais" sheet) with the 17:05 snapshot
# (previous snapshot was 16:10) and re-ranks a few countries whose updated
# totals change their sort order in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = 'Datos actualizados a 30 de Mayo de 2020 a las 17:05'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1797032
$ws.Range("C4").Value = 3502
$ws.Range("D4").Value = 519717
$ws.Range("E4").Value = 1172682
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 104633

# Alemania (row 11)
$ws.Range("B11").Value = 183118
$ws.Range("C11").Value = 99
$ws.Range("E11").Value = 9620

# India (row 12)
$ws.Range("B12").Value = 176370
$ws.Range("C12").Value = 2879
$ws.Range("D12").Value = 84500
$ws.Range("E12").Value = 86846
$ws.Range("G12").Value = 44
$ws.Range("H12").Value = 5024

# Singapur (row 29)
$ws.Range("D29").Value = 20727
$ws.Range("E29").Value = 13616

# Rumania (row 41)
$ws.Range("E41").Value = 4828
$ws.Range("G41").Value = 11
$ws.Range("H41").Value = 1259

# Kazajistan (row 55)
$ws.Range("E55").Value = 5287
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 38

# Moldavia (row 62)
$ws.Range("B62").Value = 8098
$ws.Range("C62").Value = 202
$ws.Range("E62").Value = 3352
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 291

# Irak (row 68)
$ws.Range("B68").Value = 6179
$ws.Range("C68").Value = 306
$ws.Range("D68").Value = 3110
$ws.Range("E68").Value = 2874
$ws.Range("G68").Value = 10
$ws.Range("H68").Value = 195

# Cuba (row 91) - updated totals
$ws.Range("B91").Value = 2025
$ws.Range("C91").Value = 20
$ws.Range("D91").Value = 1795
$ws.Range("E91").Value = 147
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 83

# Kenia moves above Estonia/Somalia/Islandia (rows 92-95 shift down)
$ws.Range("A92").Value = 'Kenia'
$ws.Range("B92").Value = 1888
$ws.Range("C92").Value = 143
$ws.Range("D92").Value = 438
$ws.Range("E92").Value = 1388
$ws.Range("H92").Value = 62

$ws.Range("A93").Value = 'Estonia'
$ws.Range("B93").Value = 1865
$ws.Range("C93").Value = 6
$ws.Range("D93").Value = 1622
$ws.Range("E93").Value = 176
$ws.Range("H93").Value = 67

$ws.Range("A94").Value = 'Somalia'
$ws.Range("B94").Value = 1828
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 310
$ws.Range("E94").Value = 1446
$ws.Range("H94").Value = 72

$ws.Range("A95").Value = 'Islandia'
$ws.Range("B95").Value = 1806
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = 1794
$ws.Range("E95").Value = 2
$ws.Range("H95").Value = 10

# Sri Lanka (row 101)
$ws.Range("B101").Value = 1566
$ws.Range("C101").Value = 8
$ws.Range("E101").Value = 775

# Republica de Chipre (row 120)
$ws.Range("B120").Value = 944
$ws.Range("C120").Value = 2
$ws.Range("E120").Value = 143

# Reunion (row 138)
$ws.Range("B138").Value = 471
$ws.Range("C138").Value = 1
$ws.Range("E138").Value = 59

# Cabo Verde (row 145)
$ws.Range("B145").Value = 421
$ws.Range("C145").Value = 16
$ws.Range("D145").Value = 167
$ws.Range("E145").Value = 250

# Fiyi / Curazao swap (rows 198-199)
$ws.Range("A198").Value = 'Fiyi'
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = 'Curazao'
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# Seychelles / Montserrat swap (rows 210-211)
$ws.Range("A210").Value = 'Seychelles'
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = 'Montserrat'
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Papua Nueva Guinea / Islas Virgenes Britanicas swap (rows 213-214)
$ws.Range("A213").Value = 'Papua Nueva Guinea'
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = 'Islas Virgenes Britanicas'
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
